$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 133-135 previously had stray empty-string cells in C/D (leftover
# blank form-answer columns). Drop them - only A (name) and B (presence)
# remain populated for these rows.
$ws.Range("C133").ClearContents()
$ws.Range("D133").ClearContents()
$ws.Range("C134").ClearContents()
$ws.Range("D134").ClearContents()
$ws.Range("C135").ClearContents()
$ws.Range("D135").ClearContents()

# Append the new index/response rows 136-141 (name + presence only)
$ws.Range("A136").Value = "богдана"
$ws.Range("B136").Value = "Відсутній"

$ws.Range("A137").Value = "лох"
$ws.Range("B137").Value = "Відсутній"

$ws.Range("A138").Value = "лох"
$ws.Range("B138").Value = "Відсутній"

$ws.Range("A139").Value = "лох"
$ws.Range("B139").Value = "Відсутній"

$ws.Range("A140").Value = "лохушка"
$ws.Range("B140").Value = "Відсутній"

$ws.Range("A141").Value = "iiandjdmd"
$ws.Range("B141").Value = "Відсутній"

# Row 142 mirrors the older rows: name + presence, plus the same trailing
# empty-string C/D cells that rows 133-135 used to have before being
# trimmed above.
$ws.Range("A142").Value = "iiandjdmd"
$ws.Range("B142").Value = "Відсутній"
$ws.Range("C142").Value = "'"
$ws.Range("C142").ClearFormats()
$ws.Range("D142").Value = "'"
$ws.Range("D142").ClearFormats()
